$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 (A-E)
$ws.Range("A6").Value = 43521
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = "Edwin"
$ws.Range("C6").Value = "App requests"
$ws.Range("D6").Value = "Interview"
$ws.Range("E6").Value = "2 hours (05:30 - 07:00)"

# New row 7 (A-E)
$ws.Range("A7").Value = 26
$ws.Range("B7").Value = "Edwin"
$ws.Range("C7").Value = "App requests"
$ws.Range("D7").Value = "Requeriments"
$ws.Range("E7").Value = "1.5 hours (3:30 - 05:00)"

# Column F filled in afterwards for both rows
$ws.Range("F6").Value = "Nice, I made a small talk with different Stakeholders"
$ws.Range("F7").Value = "Ok, I made progress in the requirements"

$ws.Range("C6").Select()
